# Applies the "output generated at 456a3b4" gh-pages update:
# a new exhibition ("南宁·炸裂次元动漫嘉年华") is inserted as the 5th data row
# (2024-06-01) on both the "展览" and "全部类型" sheets, pushing the later
# rows down by one, and several "想去人数" (F column) counters are bumped.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    # Force the cell to stay plain text even when the string looks like a
    # date (e.g. "2024-06-01") so Excel doesn't silently coerce it into a
    # date serial number. ClearFormats() afterwards drops the temporary
    # "@" number format so no stray style lingers on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Insert-ExhibitionRow {
    param($ws, [int]$lastRow)

    # Push rows 5..lastRow down to 6..(lastRow+1) and open up row 5.
    $ws.Rows.Item(5).Insert()

    # Give the new row's index cell (column A) the same look (bold, bordered,
    # centered) as the rest of column A by copying the formatting from the
    # row above; the real index value is written further down.
    $ws.Cells.Item(4,1).Copy($ws.Cells.Item(5,1))

    # Populate the newly-opened row with the new exhibition's data.
    Set-TextValue $ws.Cells.Item(5,2) "2024-06-01"
    $ws.Cells.Item(5,3).Value = "南宁·炸裂次元动漫嘉年华"
    $ws.Cells.Item(5,4).Value = "星光大道4号(南宁剧场地铁站D口步行220米) 文创大厦"
    $ws.Cells.Item(5,5).Value = "2024.06.01 10:00-06.02 17:00"
    $ws.Cells.Item(5,6).Value = 3
    $ws.Cells.Item(5,7).Value = 45
    $ws.Cells.Item(5,8).Value = "https://show.bilibili.com/platform/detail.html?id=85674"
    $ws.Cells.Item(5,9).Value = "//i2.hdslb.com/bfs/openplatform/202405/nYLsFLaz1715339607741.jpeg"

    # Renumber the index column (A) for the new row and every row that
    # shifted down so it stays a contiguous 1.. sequence.
    for ($r = 5; $r -le ($lastRow + 1); $r++) {
        $ws.Cells.Item($r,1).Value = $r - 1
    }
}

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) — had rows 2..12, gains the new row at 5,
# becoming rows 2..13.
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
Insert-ExhibitionRow $wsExpo 12

$wsExpo.Cells.Item(4,6).Value = 264    # 南宁·第五人格Only1.0           263 -> 264
$wsExpo.Cells.Item(6,6).Value = 3127   # 南宁·AP动漫游戏嘉年华          3117 -> 3127
$wsExpo.Cells.Item(7,6).Value = 2070   # 南宁·布谷鸟动漫展4th           2068 -> 2070
$wsExpo.Cells.Item(8,6).Value = 398    # 南宁·恋与深空only              397 -> 398
$wsExpo.Cells.Item(9,6).Value = 146    # 南宁·小蜜蜂动漫嘉年华2.0       145 -> 146
$wsExpo.Cells.Item(10,6).Value = 1170  # 南宁·AB动漫游戏嘉年华          1163 -> 1170
$wsExpo.Cells.Item(12,6).Value = 982   # 良牙夏典                       953 -> 982
$wsExpo.Cells.Item(13,6).Value = 83    # 南宁·蔚蓝档案only              81 -> 83

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types) — had rows 2..13, gains the new row at 5,
# becoming rows 2..14.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Insert-ExhibitionRow $wsAll 13

$wsAll.Cells.Item(4,6).Value = 264    # 南宁·第五人格Only1.0           263 -> 264
$wsAll.Cells.Item(6,6).Value = 3127   # 南宁·AP动漫游戏嘉年华          3117 -> 3127
$wsAll.Cells.Item(7,6).Value = 2070   # 南宁·布谷鸟动漫展4th           2068 -> 2070
$wsAll.Cells.Item(8,6).Value = 398    # 南宁·恋与深空only              397 -> 398
$wsAll.Cells.Item(10,6).Value = 146   # 南宁·小蜜蜂动漫嘉年华2.0       145 -> 146
$wsAll.Cells.Item(11,6).Value = 1170  # 南宁·AB动漫游戏嘉年华          1163 -> 1170
$wsAll.Cells.Item(13,6).Value = 982   # 良牙夏典                       953 -> 982
$wsAll.Cells.Item(14,6).Value = 83    # 南宁·蔚蓝档案only              81 -> 83
